# PIN-847: fixed tests to work with short version of locations
#
# The import-test fixture used 4 distinct admin levels (Adm1..Adm4):
#   Aleppo / Jebel Saman / Daret Azza / Kafrantin
# It is changed to use a "short" location hierarchy instead, where
# Adm1-Adm3 collapse to the same value and only Adm4 differs:
#   Al-Hasakeh / Al-Hasakeh / Al-Hasakeh / Al Berij

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$adm1 = "Al-Hasakeh"
$adm2 = "Al-Hasakeh"
$adm3 = "Al-Hasakeh"
$adm4 = "Al Berij"

for ($row = 3; $row -le 7; $row++) {
    $ws.Range("M$row").Value = $adm1
    $ws.Range("N$row").Value = $adm2
    $ws.Range("O$row").Value = $adm3
    $ws.Range("P$row").Value = $adm4
}
